# Generate Report for Handback
#
# For each language sheet (zh-cn, de-de):
#   - Status (col C) moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" for every data row.
#   - Two new columns get populated for every data row:
#       F = Latest Target File    (same file as the source .md in col A)
#       G = Latest Handback File  (same file as the handoff .xlf in col D)
#     both rendered as hyperlinks, same as their column-A / column-D siblings.
#   - Latest Handback DateTime (col H) moves from the placeholder
#     "0001-01-01 00:00:00" to a real timestamp.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

function Get-HyperlinkAddress($ws, $cell) {
    $target = $cell.Address()
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $target) {
            return $h.Address
        }
    }
    return $null
}

function Update-LangSheet($SheetName, $HandbackDateTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    for ($row = 2; $row -le 3; $row++) {

        # Status column
        $ws.Cells.Item($row, 3).Value = $newStatus

        # Source md file name/url (column A) -> mirror into Latest Target File (F)
        $srcCell = $ws.Cells.Item($row, 1)
        $srcName = $srcCell.Value2
        $srcUrl = Get-HyperlinkAddress $ws $srcCell
        $ws.Cells.Item($row, 6).Value = $srcName
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $srcUrl, "", "", $srcName) | Out-Null

        # Handoff xlf file name/url (column D) -> mirror into Latest Handback File (G)
        $targetCell = $ws.Cells.Item($row, 4)
        $targetName = $targetCell.Value2
        $targetUrl = Get-HyperlinkAddress $ws $targetCell
        $ws.Cells.Item($row, 7).Value = $targetName
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 7), $targetUrl, "", "", $targetName) | Out-Null

        # Latest Handback DateTime
        $ws.Cells.Item($row, 8).Value = $HandbackDateTime
    }
}

Update-LangSheet "zh-cn" "2016-03-20 00:36:26"
Update-LangSheet "de-de" "2016-03-20 00:36:32"

# The Overview sheet mirrors each language's Status in columns B (zh-cn) and
# C (de-de) for every file row -- keep them in sync with the per-language
# sheets above.
$overview = $wb.Worksheets.Item("Overview")
for ($row = 2; $row -le 3; $row++) {
    $overview.Cells.Item($row, 2).Value = $newStatus
    $overview.Cells.Item($row, 3).Value = $newStatus
}
